$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "580÷3=193, 1"
$t.Cell(1, 2).Range.Text = "384÷2=192, 0"
$t.Cell(1, 3).Range.Text = "960÷9=106, 6"
$t.Cell(1, 4).Range.Text = "365÷8=45, 5"
$t.Cell(1, 5).Range.Text = "767÷9=85, 2"
$t.Cell(5, 1).Range.Text = "926÷6=154, 2"
$t.Cell(5, 2).Range.Text = "770÷4=192, 2"
$t.Cell(5, 3).Range.Text = "467÷6=77, 5"
$t.Cell(5, 4).Range.Text = "962÷2=481, 0"
$t.Cell(5, 5).Range.Text = "395÷6=65, 5"
$t.Cell(9, 1).Range.Text = "574÷8=71, 6"
$t.Cell(9, 2).Range.Text = "594÷6=99, 0"
$t.Cell(9, 3).Range.Text = "297÷9=33, 0"
$t.Cell(9, 4).Range.Text = "897÷4=224, 1"
$t.Cell(9, 5).Range.Text = "164÷5=32, 4"
$t.Cell(13, 1).Range.Text = "466÷2=233, 0"
$t.Cell(13, 2).Range.Text = "671÷8=83, 7"
$t.Cell(13, 3).Range.Text = "225÷8=28, 1"
$t.Cell(13, 4).Range.Text = "360÷8=45, 0"
$t.Cell(13, 5).Range.Text = "387÷4=96, 3"
$t.Cell(17, 1).Range.Text = "356÷4=89, 0"
$t.Cell(17, 2).Range.Text = "573÷6=95, 3"
$t.Cell(17, 3).Range.Text = "278÷5=55, 3"
$t.Cell(17, 4).Range.Text = "192÷8=24, 0"
$t.Cell(17, 5).Range.Text = "866÷5=173, 1"

Write-Output "done"
